$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 content updates:
#  - C2 previously held the stray value "42" -> corrected to the issue type "Feature"
#  - E2 previously held "Hombergs, Tom" -> corrected to "hombergs"
$ws.Range("C2").Value = "Feature"
$ws.Range("E2").Value = "hombergs"

# Update the active selection/view: unpin the frozen top-left scroll position
# and move the active cell/selection from J4 to J2.
$ws.Range("J2").Select()
